$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-label the existing header cells -----------------------------
# Two new tax-rate column pairs ("%10" and "%20") are inserted into the
# KDV/Matrah header sequence. Rather than hand-edit the shared-string
# table, just retype every header cell from H1 onward with its new
# final caption (exactly what a user fixing the report in Excel would
# do); Excel manages the shared strings internally.
$ws.Range("H1").Value2 = "Matrah 10"
$ws.Range("I1").Value2 = "KDV10"
$ws.Range("J1").Value2 = "Matrah 18"
$ws.Range("K1").Value2 = "KDV 18"
$ws.Range("L1").Value2 = "Matrah 20"
$ws.Range("M1").Value2 = "KDV 20"
$ws.Range("N1").Value2 = "Matrah"
$ws.Range("O1").Value2 = "KDV"
$ws.Range("P1").Value2 = "KDV Oranı"
$ws.Range("Q1").Value2 = "Hesap Kodu"
$ws.Range("R1").Value2 = "Doviz"
$ws.Range("S1").Value2 = "Doviz Matrah 8"
$ws.Range("T1").Value2 = "Doviz KDV8"
$ws.Range("U1").Value2 = "Doviz Matrah 18"
$ws.Range("V1").Value2 = "Doviz KDV18"
$ws.Range("W1").Value2 = "Doviz Matrah"
$ws.Range("X1").Value2 = "Doviz KDV"

# --- 2) Re-apply the correct header cell formatting ---------------------
# Copy formats from cells that already carry the desired look so no new
# style entries are created and the existing style palette is reused.

# Q1 & R1 should use the "text" header style that M1/N1 used to have.
$ws.Range("M1").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null

# U1:X1 should use the plain thin-border header style that S1 (and
# previously O1/Q1/R1/T1) uses.
$ws.Range("S1").Copy() | Out-Null
$ws.Range("U1:X1").PasteSpecial(-4122) | Out-Null

# T1 should use the bold/black header style that P1 used to have.
$ws.Range("P1").Copy() | Out-Null
$ws.Range("T1").PasteSpecial(-4122) | Out-Null

# M1, N1, O1, P1 should use the plain hairline header style that A1 uses.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("M1:P1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- 3) Tidy up leftover clipboard state / selection ---------------------
$ws.Range("J5").Select() | Out-Null

# --- 4) Column widths (best-effort; Excel snaps these to whole pixels) --
$widths = @{
    1  = 5.46
    2  = 9.47
    3  = 4.36
    4  = 7.82
    5  = 8.38
    6  = 8.66
    7  = 6.29
    8  = 9.63
    9  = 9.63
    10 = 9.63
    11 = 7.82
    12 = 9.63
    13 = 7.82
    14 = 7.13
    15 = 5.32
    16 = 10.32
    17 = 11.43
    18 = 6.16
    19 = 13.82
    20 = 11.43
    21 = 14.77
    22 = 12.41
    23 = 12.27
    24 = 10.46
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}
